$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header row (row 3) - shared string labels
$ws.Range("B3").Value = "history10"
$ws.Range("C3").Value = "history30"
$ws.Range("D3").Value = "history50"
$ws.Range("E3").Value = "history70"
$ws.Range("F3").Value = "history90"

# Row 5: images/s
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 70
$ws.Range("F5").Value = 90

# Row 6: recuperation profondeur
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 9

# Row 7: fiabilite height
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

# Row 8: proportion de cibles
$ws.Range("B8").Value = 0.0509259259259259
$ws.Range("C8").Value = 0.0509259259259259
$ws.Range("D8").Value = 0.0509259259259259
$ws.Range("E8").Value = 0.0509259259259259
$ws.Range("F8").Value = 0.0509259259259259

# Row 9: nombre de points
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 9

# Row 10: nombre de pertes
$ws.Range("B10").Value = 43.2
$ws.Range("C10").Value = 43.2
$ws.Range("D10").Value = 43.2
$ws.Range("E10").Value = 43.2
$ws.Range("F10").Value = 43.2

# Row 11: indice de performance
$ws.Range("B11").Value = -499
$ws.Range("C11").Value = -499
$ws.Range("D11").Value = -499
$ws.Range("E11").Value = -499
$ws.Range("F11").Value = -499

# Row 12: extra zero row (not referenced by any chart series)
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

$wb.Save()
